$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TreatmentInstruction")

# Three ZIB concepts (rows 3, 4, 5 in the "Notes/Issues" mapping table) are
# now made into extensions: update the "Maps to" (J) and "Notes/Issues" (K)
# columns accordingly.
$ws.Range("J3").Value = "Consent.extension"
$ws.Range("K3").Value = "gForge #13540"

$ws.Range("J4").Value = "Consent.extension"
$ws.Range("K4").Value = "gForge #13540"

$ws.Range("J5").Value = "Concent.except.extenstion"
$ws.Range("K5").Value = "gForge #13540"

# Move the active selection/cursor to H22, as last left by the author.
$ws.Range("H22").Select()
